$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $origFormat = $Range.NumberFormat
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.NumberFormat = $origFormat
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '52.197.43'
Set-TextValue $ws.Range('E2') '  +0.03%  '
Set-TextValue $ws.Range('D3') '2.800.09'
Set-TextValue $ws.Range('E3') '  -0.97%  '
Set-TextValue $ws.Range('E4') '  +0.14%  '
Set-TextValue $ws.Range('D5') '362.79'
Set-TextValue $ws.Range('E5') '  +1.62%  '
Set-TextValue $ws.Range('D6') '110.20'
Set-TextValue $ws.Range('E6') '  -3.64%  '
Set-TextValue $ws.Range('E7') '  +2.84%  '
Set-TextValue $ws.Range('D8') '1.00'
Set-TextValue $ws.Range('E8') '  +0.12%  '
Set-TextValue $ws.Range('D9') '0.598'
Set-TextValue $ws.Range('E9') '  -0.97%  '
Set-TextValue $ws.Range('D10') '40.41'
Set-TextValue $ws.Range('E10') '  -3.67%  '
Set-TextValue $ws.Range('D11') '0.0855'
Set-TextValue $ws.Range('E11') '  +0.85%  '
Set-TextValue $ws.Range('E12') '  +0.31%  '
Set-TextValue $ws.Range('D13') '19.64'
Set-TextValue $ws.Range('E13') '  -1.33%  '
Set-TextValue $ws.Range('D14') '7.69'
Set-TextValue $ws.Range('E14') '  -1.57%  '
Set-TextValue $ws.Range('D15') '3.238.47'
Set-TextValue $ws.Range('E15') '  -0.61%  '
Set-TextValue $ws.Range('D16') '2.813.10'
Set-TextValue $ws.Range('E16') '  -0.73%  '
Set-TextValue $ws.Range('D17') '0.936'
Set-TextValue $ws.Range('E17') '  +4.45%  '
Set-TextValue $ws.Range('D18') '52.090.43'
Set-TextValue $ws.Range('E18') '  +0.11%  '
Set-TextValue $ws.Range('D19') '7.41'
Set-TextValue $ws.Range('E19') '  +1.73%  '
Set-TextValue $ws.Range('D20') '3.14'
Set-TextValue $ws.Range('E20') '  -1.26%  '
Set-TextValue $ws.Range('D21') '13.22'
Set-TextValue $ws.Range('E21') '  -3.33%  '
Set-TextValue $ws.Range('E22') '  -0.12%  '
Set-TextValue $ws.Range('D23') '273.97'
Set-TextValue $ws.Range('E23') '  +1.42%  '
Set-TextValue $ws.Range('D24') '69.91'
Set-TextValue $ws.Range('E24') '  +0.44%  '
Set-TextValue $ws.Range('E25') '  -0.44%  '
Set-TextValue $ws.Range('D26') '26.79'
Set-TextValue $ws.Range('E26') '  -0.18%  '
Set-TextValue $ws.Range('E27') '  -0.07%  '
Set-TextValue $ws.Range('D28') '10.24'
Set-TextValue $ws.Range('E28') '  +0.04%  '
Set-TextValue $ws.Range('E29') '  -0.77%  '
Set-TextValue $ws.Range('E30') '  +2.11%  '
Set-TextValue $ws.Range('D31') '0.0475'
Set-TextValue $ws.Range('E31') '  +8.83%  '
Set-TextValue $ws.Range('D32') '51.62'
Set-TextValue $ws.Range('E32') '  +1.85%  '
Set-TextValue $ws.Range('D33') '34.59'
Set-TextValue $ws.Range('E33') '  +2.13%  '
Set-TextValue $ws.Range('D34') '5.79'
Set-TextValue $ws.Range('E34') '  -1.03%  '
Set-TextValue $ws.Range('D35') '5.47'
Set-TextValue $ws.Range('E35') '  +11.59%  '
Set-TextValue $ws.Range('D36') '0.0844'
Set-TextValue $ws.Range('E36') '  +1.87%  '
Set-TextValue $ws.Range('D37') '1.00'
Set-TextValue $ws.Range('E37') '  +0.20%  '
Set-TextValue $ws.Range('D38') '3.22'
Set-TextValue $ws.Range('E38') '  +0.63%  '
Set-TextValue $ws.Range('D39') '18.41'
Set-TextValue $ws.Range('E39') '  -0.42%  '
Set-TextValue $ws.Range('D40') '2.02'
Set-TextValue $ws.Range('E40') '  -3.98%  '
Set-TextValue $ws.Range('D41') '2.59'
Set-TextValue $ws.Range('E41') '  +0.61%  '
Set-TextValue $ws.Range('E42') '  +0.02%  '
Set-TextValue $ws.Range('D43') '125.03'
Set-TextValue $ws.Range('E43') '  -1.31%  '
Set-TextValue $ws.Range('D44') '2.25'
Set-TextValue $ws.Range('E44') '  -1.96%  '
Set-TextValue $ws.Range('D45') '22.07'
Set-TextValue $ws.Range('E45') '  -5.47%  '
Set-TextValue $ws.Range('D46') '2.071.82'
Set-TextValue $ws.Range('E46') '  +1.31%  '
Set-TextValue $ws.Range('D47') '3.28'
Set-TextValue $ws.Range('E47') '  -2.03%  '
Set-TextValue $ws.Range('E48') '  +0.09%  '
Set-TextValue $ws.Range('D49') '5.77'
Set-TextValue $ws.Range('E49') '  +1.33%  '
Set-TextValue $ws.Range('D50') '0.948'
Set-TextValue $ws.Range('E50') '  -0.28%  '
Set-TextValue $ws.Range('D51') '9.08'
Set-TextValue $ws.Range('E51') '  +2.21%  '
